$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every value (even dates/numbers) as literal text (t="str").
# Force the new row's cells to Text format so Excel keeps them as text
# instead of converting to real dates/numbers.
$ws.Range("A9:E9").NumberFormat = "@"

$ws.Range("A9").Value = "2025-09-25"
$ws.Range("B9").Value = "Pick 4"
$ws.Range("C9").Value = "250925"
$ws.Range("D9").Value = "5-2-5-5"
$ws.Range("E9").Value = "2025-09-25T21:37:53.171+04:00"
